$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Restructure the script flow ---------------------------------------

# Insert two new rows (10:11) for a pair of new "swipe" steps that now
# happen right after the team walks to location2.
$ws.Rows("10:11").Insert()

# Copy the formatting of the row above (location2, row 9) onto the two new
# rows so the borders / alignment match the rest of the table.
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new swipe rows.
$ws.Range("A10").Value2 = "500, 130, 500, 450"
$ws.Range("B10").Value2 = "swipe"
$ws.Range("A11").Value2 = "500, 130, 500, 450"
$ws.Range("B11").Value2 = "swipe"

# The old stand-alone "swipe" row (originally row 16, now pushed to row 18
# by the insert above) is no longer needed - that swipe now happens earlier.
$ws.Rows("18:18").Delete()

# --- 2. Correct / update coordinates and values ----------------------------

$ws.Range("A9").Value2  = "376, 200"    # location2 coordinate fix
$ws.Range("A12").Value2 = "400, 422"    # enemy2 coordinate fix
$ws.Range("A15").Value2 = "400, 422"    # location3 coordinate fix
$ws.Range("A16").Value2 = "502, 266"    # blank coordinate fix
$ws.Range("D17").Value2 = 15            # increase wait_seconds timeout
$ws.Range("A18").Value2 = "502, 266"    # location4 coordinate fix

# --- 3. View / selection state ---------------------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D17").Select()
